$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 0. Big italic paragraph: split into runs around "ideas" / "and also" / "has to" ---
$p = $d.Paragraphs(5)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $W + '>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Upon signing the employment contract with Capgemini/Matiq, Consultant has agreed to a confidentiality clause, where he/she has agreed to undertake professional secrecy and not to make unauthorized use of information about Capgemini/Matiq, its clients, work methods, property, </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ideas</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> or anything else pertaining to Capgemini/Matiq of which he/she has received knowledge as of the date of signing the employment agreement and during the whole period of employment. Professional secrecy applies to all information that can be deemed sensitive or constitute a business secret, the disclosure or use of which can cause damage to clients, personnel or to Capgemini/Matiq. Professional secrecy applies for the term of employment </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>and also</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> after the employment has ceased. Consultant </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>has to</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> carry out his/her work duties either at Capgemini’s or the client’s premises, and for this reason work material may not be kept at any other location.</w:t></w:r>' +
  '</w:p>'
$rng.InsertXML($xml)

# --- 1. "Rajneesh Hajela " paragraph -> "____________ " ---
$p = $d.Paragraphs(9)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $W + '><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">____________ </w:t></w:r></w:p>'
$rng.InsertXML($xml)

# --- 2. "This confirmation does not..." paragraph: split around "i.e." ---
$p = $d.Paragraphs(14)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $W + '>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">This confirmation does not establish a legal basis for personal liabilities against the consultant, </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>i.e.</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> any complaints, claims, liabilities may only be directed against Capgemini according to the applicable MSA. The consultant has liabilities according to work contract with Capgemini.</w:t></w:r>' +
  '</w:p>'
$rng.InsertXML($xml)

# --- 3. "Mumbai 23/July/2024" + signature-line paragraphs ---
# Replace paragraph 18 ("Mumbai...") with an empty paragraph, and paragraph 19
# (the "____________________ ... Rajneesh" line) with new underscore placeholders.
$p18 = $d.Paragraphs(18)
$rng18 = $d.Range($p18.Range.Start, $p18.Range.End - 1)
$rng18.InsertXML('<w:p ' + $W + '/>')

$p19 = $d.Paragraphs(19)
$rng19full = $d.Range($p19.Range.Start, $p19.Range.End)
$xml19 = '<w:p ' + $W + '>' +
  '<w:r><w:t xml:space="preserve">________________________________ </w:t></w:r>' +
  '<w:r><w:tab/></w:r>' +
  '<w:r><w:tab/><w:t>________________________________</w:t></w:r>' +
  '</w:p>'
$rng19full.InsertXML($xml19)
